# Update "想去人数" (F column) counts on both the "展览" sheet and the
# aggregated "全部类型" sheet to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 229
$wsExhibit.Range("F4").Value = 844
$wsExhibit.Range("F5").Value = 74
$wsExhibit.Range("F6").Value = 31

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 229
$wsAll.Range("F5").Value = 844
$wsAll.Range("F6").Value = 74
$wsAll.Range("F7").Value = 31
